# The source workbook gained one new daily price-observation row for
# "Feria Lagunitas de Puerto Montt" / Perejil. Excel's row-insert at 43
# shifts the existing rows 43-167 down to 44-168 (preserving their data
# untouched) and we fill the freshly inserted row 43 with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Range('A43').Value = 4
$ws.Range('B43').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C43').Value = 'Los Lagos'
$ws.Range('D43').Value = 44497
$ws.Range('E43').Value = 10
$ws.Range('F43').Value = 100112044
$ws.Range('G43').Value = 'Perejil'
$ws.Range('H43').Value = 'Sin especificar'
$ws.Range('I43').Value = 'Primera'
$ws.Range('J43').Value = 80
$ws.Range('K43').Value = 5000
$ws.Range('L43').Value = 5000
$ws.Range('M43').Value = 5000
$ws.Range('N43').Value = '$/docena de atados (3 kilos)'
$ws.Range('O43').Value = 'Región Metropolitana'
$ws.Range('P43').Value = 1667
$ws.Range('Q43').Value = 3
$ws.Range('R43').Value = 'Hortaliza'
